$wb = $excel.ActiveWorkbook

# Column letter helper via fixed tables matching the 5 repeating blocks
$blocks = @(
  @{ Label = 1;  Value = 2;  Days = @(3,4,5,6,7,8) },    # A-H   przedmiot
  @{ Label = 10; Value = 11; Days = @(12,13,14,15,16,17) }, # J-Q   grupa
  @{ Label = 19; Value = 20; Days = @(21,22,23,24,25,26) }, # S-Z   nauczyciel
  @{ Label = 28; Value = 29; Days = @(30,31,32,33,34,35) }, # AB-AI klasa
  @{ Label = 37; Value = 38; Days = @(39,40,41,42,43,44) }  # AK-AR sala
)

foreach ($ws in $wb.Worksheets) {
  $maxRow = $ws.UsedRange.Rows.Count()

  foreach ($blk in $blocks) {
    $valueCol = $blk.Value
    $labelCol = $blk.Label
    $days = $blk.Days

    $count = 0
    $sums = @(0,0,0,0,0,0)
    $lastRow = 1

    for ($r = 2; $r -le $maxRow; $r++) {
      $v = $ws.Cells.Item($r, $valueCol).Value()
      if ($v -ne $null) {
        $count = $count + 1
        $lastRow = $r
        for ($i = 0; $i -lt 6; $i++) {
          $dv = $ws.Cells.Item($r, $days[$i]).Value()
          if ($dv -ne $null) {
            $sums[$i] = $sums[$i] + $dv
          }
        }
      }
    }

    $razemRow = $lastRow + 1

    # Label cell ("Razem") - styled like the other header/label cells
    $labelCell = $ws.Cells.Item($razemRow, $labelCol)
    $labelCell.Value = "Razem"

    # Value cell (count) - also styled
    $valCell = $ws.Cells.Item($razemRow, $valueCol)
    $valCell.Value = $count

    # copy formatting from row-1 header cell in the same column-block (bold/border/center style)
    $ws.Cells.Item(1, $valueCol + 6).Copy()
    $labelCell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
    $ws.Cells.Item(1, $valueCol + 6).Copy()
    $valCell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

    for ($i = 0; $i -lt 6; $i++) {
      $ws.Cells.Item($razemRow, $days[$i]).Value = $sums[$i]
    }
  }
}
